$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update billed character counts (D3, D4)
$ws.Range("D3").Value = 1040
$ws.Range("D4").Value = 1039

# Move the active cell selection from D5 to D4
$ws.Range("D4").Select()
